$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename Sheet2 -> ProductTestData, then insert the remaining new sheets
#    in order right after it so that sheet2.xml..sheet10.xml map onto
#    ProductTestData..FlipkartProfileInfoTestData.
# ---------------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item(1)

$productSheet = $wb.Worksheets.Item(2)
$productSheet.Name = "ProductTestData"

$wishlistSheet = $wb.Worksheets.Add($null, $productSheet)
$wishlistSheet.Name = "WishlistTestData"

$addressSheet = $wb.Worksheets.Add($null, $wishlistSheet)
$addressSheet.Name = "AddressTestData"

$loginSheet = $wb.Worksheets.Add($null, $addressSheet)
$loginSheet.Name = "LoginTestData"

$homeSheet = $wb.Worksheets.Add($null, $loginSheet)
$homeSheet.Name = "HomePageTestData"

$logoutSheet = $wb.Worksheets.Add($null, $homeSheet)
$logoutSheet.Name = "LogoutTestData"

$cartSheet = $wb.Worksheets.Add($null, $logoutSheet)
$cartSheet.Name = "CartTestData"

$productDetailsSheet = $wb.Worksheets.Add($null, $cartSheet)
$productDetailsSheet.Name = "ProductDetailsTestData"

$flipkartProfileSheet = $wb.Worksheets.Add($null, $productDetailsSheet)
$flipkartProfileSheet.Name = "FlipkartProfileInfoTestData"

# ---------------------------------------------------------------------------
# 2. ProductTestData content
# ---------------------------------------------------------------------------
$productSheet.Range("A1").Value = "TestCaseName"
$productSheet.Range("B1").Value = "Execution Required"
$productSheet.Range("C1").Value = "keyword"
$productSheet.Range("D1").Value = "brand"
$productSheet.Range("A1:D1").Font.Bold = $true

$productSheet.Range("A2").Value = "relevantProducts"
$productSheet.Range("B2").Value = "Yes "
$productSheet.Range("C2").Value = "iphone"
$productSheet.Range("D2").Value = "APPLE"

$productSheet.Range("A3").Value = "sortProductByPriceHtoL"
$productSheet.Range("B3").Value = "Yes "
$productSheet.Range("C3").Value = "iphone"

$productSheet.Range("A4").Value = "sortProductByPriceLtoH"
$productSheet.Range("B4").Value = "Yes "
$productSheet.Range("C4").Value = "iphone"

$productSheet.Range("A5").Value = "chooseProductBrand"
$productSheet.Range("B5").Value = "Yes "
$productSheet.Range("C5").Value = "mobile"
$productSheet.Range("D5").Value = "apple"

$productSheet.Range("A6").Value = "chooseProductRating"
$productSheet.Range("B6").Value = "Yes "
$productSheet.Range("C6").Value = "iphone"

$productSheet.Columns.Item(1).ColumnWidth = 25
$productSheet.Columns.Item(2).ColumnWidth = 18.7265625
$productSheet.Columns.Item(3).ColumnWidth = 24.08984375

# ---------------------------------------------------------------------------
# 3. WishlistTestData content
# ---------------------------------------------------------------------------
$wishlistSheet.Range("A1").Value = "TestCaseName"
$wishlistSheet.Range("B1").Value = "Execution Required"
$wishlistSheet.Range("C1").Value = "keyword"
$wishlistSheet.Range("D1").Value = "brand"
$wishlistSheet.Range("E1").Value = "mob"
$wishlistSheet.Range("F1").Value = "pwd"
$wishlistSheet.Range("G1").Value = "keyword"
$wishlistSheet.Range("A1:G1").Font.Bold = $true

$wishlistSheet.Range("A2").Value = "addTowishlist"
$wishlistSheet.Range("B2").Value = "No"
$wishlistSheet.Range("C2").Value = "iphone"
$wishlistSheet.Range("D2").Value = "APPLE"
$wishlistSheet.Range("E2").Value = """8708185463"""
$wishlistSheet.Range("F2").Value = "pulkit9017"
$wishlistSheet.Range("G2").Value = "iPhone"

$wishlistSheet.Range("A3").Value = "removeFromwishlist"
$wishlistSheet.Range("B3").Value = "Yes"
$wishlistSheet.Range("C3").Value = "iphone"
$wishlistSheet.Range("E3").Value = """8708185463"""
$wishlistSheet.Range("F3").Value = "pulkit9017"

$wishlistSheet.Range("A4").Value = "getProductsInwishlist"
$wishlistSheet.Range("B4").Value = "Yes"
$wishlistSheet.Range("C4").Value = "iphone"
$wishlistSheet.Range("E4").Value = """8708185463"""
$wishlistSheet.Range("F4").Value = "pulkit9017"

$wishlistSheet.Columns.Item(1).ColumnWidth = 22.54296875
$wishlistSheet.Columns.Item(2).ColumnWidth = 23.1796875
$wishlistSheet.Columns.Item(5).ColumnWidth = 10.81640625
$wishlistSheet.Columns.Item(6).ColumnWidth = 13.81640625

# ---------------------------------------------------------------------------
# 4. Sheet selections / active cells
# ---------------------------------------------------------------------------
$sheet1.Range("C2").Select()
$productSheet.Range("A1:D3").Select()
$productDetailsSheet.Range("K17").Select()
$flipkartProfileSheet.Range("H14").Select()

# WishlistTestData is the tab that should end up active/selected.
$wishlistSheet.Range("B4:G4").Select()

$wb.Save()
